{"js": "// Insert two new bulleted list items (\"- Faire un push sur le repository\n// distant\" and \"- Faire un merge pour voir\") right after the existing\n// last list item (\"Modifier les fichiers et faire un commit\"). Using\n// InsertLocation.after on that paragraph makes the new paragraphs\n// inherit its paragraph formatting (ListParagraph style + the same\n// bullet-list numPr), matching the target structure.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Modifier les fichiers et faire un commit\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  // Fallback: the last paragraph in the document body.\n  anchor = paragraphs.items[paragraphs.items.length - 1];\n}\n\nconst first = anchor.insertParagraph(\n  \"- Faire un push sur le repository distant\",\n  Word.InsertLocation.after\n);\nfirst.insertParagraph(\"- Faire un merge pour voir\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert two new bulleted list items (\"- Faire un push sur le repository\n# distant\" and \"- Faire un merge pour voir\") right after the existing\n# last list item (\"Modifier les fichiers et faire un commit\").\n#\n# InsertParagraphAfter() on that paragraph's Range creates a new empty\n# paragraph that inherits the paragraph formatting of its neighbour\n# (ListParagraph style + the same bullet-list numPr), matching the\n# target structure exactly.\n\n$d = $word.ActiveDocument\n\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\").Trim() -eq \"Modifier les fichiers et faire un commit\") {\n        $anchor = $p\n    }\n}\nif ($anchor -eq $null) {\n    $anchor = $d.Paragraphs.Last\n}\n\n$anchor.Range.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Last\n$p1.Range.Text = \"- Faire un push sur le repository distant\"\n\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Last\n$p2.Range.Text = \"- Faire un merge pour voir\"\n"}
